# Commit "Added tables to Readme." -- refreshed calibration results for the
# "cam2" sheet (columns C:H, rows 2-31) plus the saved cell selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cam2")

# New Fx, Fy, Cx, Cy, RMS, Colinearidad values per iteration row (id 0..29)
$newValues = @(
  @(486.50099999999998, 486.17399999999998, 350.85599999999999, 178.87799999999999, 0.29776900000000001, 0.043717600000000002),
  @(486.56700000000001, 486.19900000000001, 350.77800000000002, 179.00299999999999, 0.29162199999999999, 0.034427699999999999),
  @(486.57400000000001, 486.19400000000002, 350.803, 178.99600000000001, 0.28998200000000002, 0.031148700000000001),
  @(486.65199999999999, 486.25900000000001, 350.78899999999999, 179.077, 0.28955500000000001, 0.030565700000000001),
  @(486.69299999999998, 486.298, 350.81799999999998, 179.089, 0.28958699999999998, 0.030590699999999998),
  @(486.721, 486.31799999999998, 350.79199999999997, 179.06899999999999, 0.28945199999999999, 0.030654600000000001),
  @(486.78300000000002, 486.387, 350.84199999999998, 179.001, 0.28967700000000002, 0.030684800000000002),
  @(486.86099999999999, 486.47699999999998, 350.78199999999998, 178.99700000000001, 0.28950700000000001, 0.030754400000000001),
  @(486.8, 486.41800000000001, 350.755, 178.95099999999999, 0.28951199999999999, 0.030664199999999999),
  @(486.82499999999999, 486.45, 350.81299999999999, 178.94200000000001, 0.28916199999999997, 0.0304579),
  @(486.76100000000002, 486.38799999999998, 350.82499999999999, 178.94399999999999, 0.28915200000000002, 0.030586800000000001),
  @(486.81799999999998, 486.44900000000001, 350.83800000000002, 178.99299999999999, 0.28919400000000001, 0.030442899999999998),
  @(486.84500000000003, 486.47500000000002, 350.80599999999998, 178.99100000000001, 0.288885, 0.030373000000000001),
  @(486.87, 486.49900000000002, 350.81700000000001, 178.93700000000001, 0.28920800000000002, 0.030482599999999999),
  @(486.82400000000001, 486.46300000000002, 350.79500000000002, 178.95599999999999, 0.28919600000000001, 0.030372699999999999),
  @(486.80900000000003, 486.44900000000001, 350.68900000000002, 178.947, 0.289273, 0.030582000000000002),
  @(486.87299999999999, 486.51100000000002, 350.77499999999998, 178.934, 0.28909899999999999, 0.0306945),
  @(486.81400000000002, 486.45600000000002, 350.71699999999998, 178.904, 0.28904200000000002, 0.0307149),
  @(486.88099999999997, 486.52300000000002, 350.714, 178.88900000000001, 0.289356, 0.030840800000000002),
  @(486.83, 486.46899999999999, 350.75799999999998, 178.87100000000001, 0.28943099999999999, 0.0305571),
  @(486.86700000000002, 486.50900000000001, 350.71899999999999, 178.93, 0.28964600000000001, 0.030746800000000001),
  @(486.74, 486.38200000000001, 350.66199999999998, 178.90299999999999, 0.28941800000000001, 0.0306406),
  @(486.85899999999998, 486.49599999999998, 350.68099999999998, 178.92, 0.28947000000000001, 0.030611699999999999),
  @(486.923, 486.56900000000002, 350.69900000000001, 178.89, 0.28947299999999998, 0.030481100000000001),
  @(486.923, 486.56799999999998, 350.74200000000002, 179.02, 0.289377, 0.030565800000000001),
  @(486.834, 486.47899999999998, 350.75200000000001, 178.95500000000001, 0.28930699999999998, 0.0305432),
  @(486.899, 486.541, 350.75200000000001, 178.98500000000001, 0.28930099999999997, 0.0304654),
  @(486.91800000000001, 486.565, 350.74299999999999, 178.91900000000001, 0.28923900000000002, 0.030500699999999999),
  @(486.83600000000001, 486.47300000000001, 350.80599999999998, 178.92400000000001, 0.28940300000000002, 0.0305455),
  @(486.80399999999997, 486.44900000000001, 350.71300000000002, 178.95699999999999, 0.28945100000000001, 0.030612299999999999)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
  $row = 2 + $i
  $vals = $newValues[$i]
  for ($col = 0; $col -lt $vals.Length; $col++) {
    # columns C..H => 3..8
    $ws.Cells.Item($row, 3 + $col).Value = $vals[$col]
  }
}

# Restore the active-cell selection saved with the workbook
[void]$ws.Range("J22").Select()
